$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select E8 first so the saved sheetView records it as the active cell
# (reproduces the <selection activeCell="E8" sqref="E8"/> added to sheetView)
$ws.Range("E8").Select()

# Update the greeting text for rule R10 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"
